$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "Dropped by Golem. Material used to buy Loot Boxes"
$ws.Range("B8").Value  = "Dropped by Golem. Material used to buy Premium Loot Boxes"
$ws.Range("B9").Value  = "Dropped by Panther. Material used to buy Loot Boxes"
$ws.Range("B10").Value = "Dropped by Panther. Material used to buy Loot Boxes"
$ws.Range("B11").Value = "Dropped by Panther. Material used to by Premium Loot Boxes"
$ws.Range("B12").Value = "Dropped by Treant. Material used to buy Loot Boxes"
$ws.Range("B13").Value = "Dropped by Treant. Material used to buy Premium Loot Boxes"
